# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Every value in columns D (Price) and E (Volume(1h)) is stored as plain TEXT
# in the source workbook (even cells that look numeric, e.g. "1.001"), so we
# must write them back as text too - not let Excel auto-coerce them to
# numbers - while leaving each cell's style/format untouched.
#
# Prefixing the literal with a single quote forces Excel to store it as text
# (quoted-text entry), but that also stamps the cell with a "quotePrefix"
# style. Since the original cells carry no such flag, we snapshot and restore
# the cell's Style around the write so the net style stays identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "29.986.35"
Set-TextValue $ws.Range("E2") "  -0.10%  "

Set-TextValue $ws.Range("D3") "1.892.59"
Set-TextValue $ws.Range("E3") "  -0.83%  "

Set-TextValue $ws.Range("D4") "0.9996"
Set-TextValue $ws.Range("E4") "  +0.15%  "

Set-TextValue $ws.Range("D5") "0.8310"
Set-TextValue $ws.Range("E5") "  +6.13%  "

Set-TextValue $ws.Range("D6") "242.06"
Set-TextValue $ws.Range("E6") "  +0.08%  "

Set-TextValue $ws.Range("D7") "0.9994"
Set-TextValue $ws.Range("E7") "  +0.10%  "

Set-TextValue $ws.Range("D8") "0.3249"
Set-TextValue $ws.Range("E8") "  +3.49%  "

Set-TextValue $ws.Range("D9") "26.43"
Set-TextValue $ws.Range("E9") "  +1.38%  "

Set-TextValue $ws.Range("D10") "0.07024"
Set-TextValue $ws.Range("E10") "  +1.59%  "

Set-TextValue $ws.Range("D11") "0.08050"
Set-TextValue $ws.Range("E11") "  +0.57%  "

Set-TextValue $ws.Range("D12") "0.7499"
Set-TextValue $ws.Range("E12") "  +0.12%  "

Set-TextValue $ws.Range("D13") "1.892.54"
Set-TextValue $ws.Range("E13") "  -0.77%  "

Set-TextValue $ws.Range("D14") "5.230"
Set-TextValue $ws.Range("E14") "  +0.17%  "

Set-TextValue $ws.Range("D15") "92.41"
Set-TextValue $ws.Range("E15") "  -0.17%  "

Set-TextValue $ws.Range("D16") "29.969.90"
Set-TextValue $ws.Range("E16") "  -0.15%  "

Set-TextValue $ws.Range("D17") "14.06"
Set-TextValue $ws.Range("E17") "  +0.58%  "

Set-TextValue $ws.Range("D18") "5.865"
Set-TextValue $ws.Range("E18") "  -1.47%  "

Set-TextValue $ws.Range("D19") "245.62"
Set-TextValue $ws.Range("E19") "  -1.03%  "

Set-TextValue $ws.Range("D20") "0.000007756"
Set-TextValue $ws.Range("E20") "  +0.13%  "

Set-TextValue $ws.Range("D21") "0.9996"
Set-TextValue $ws.Range("E21") "  +0.10%  "

Set-TextValue $ws.Range("D22") "2.139.21"
Set-TextValue $ws.Range("E22") "  -0.70%  "

Set-TextValue $ws.Range("D23") "0.9997"
Set-TextValue $ws.Range("E23") "  +0.12%  "

Set-TextValue $ws.Range("D24") "6.958"
Set-TextValue $ws.Range("E24") "  -0.43%  "

Set-TextValue $ws.Range("D25") "0.1589"
Set-TextValue $ws.Range("E25") "  +19.80%  "

Set-TextValue $ws.Range("D26") "165.55"
Set-TextValue $ws.Range("E26") "  -1.32%  "

Set-TextValue $ws.Range("D27") "9.208"
Set-TextValue $ws.Range("E27") "  -1.83%  "

Set-TextValue $ws.Range("D28") "18.85"
Set-TextValue $ws.Range("E28") "  -0.38%  "

Set-TextValue $ws.Range("D29") "2.096"
Set-TextValue $ws.Range("E29") "  +2.10%  "

Set-TextValue $ws.Range("D30") "1.362"
Set-TextValue $ws.Range("E30") "  -2.26%  "

# Row 31 (PancakeSwap): only the volume figure changed, price stayed "1.516".
Set-TextValue $ws.Range("E31") "  -0.14%  "

Set-TextValue $ws.Range("D32") "4.274"
Set-TextValue $ws.Range("E32") "  -0.89%  "

Set-TextValue $ws.Range("D33") "0.05713"
Set-TextValue $ws.Range("E33") "  +7.86%  "

Set-TextValue $ws.Range("D34") "4.069"
Set-TextValue $ws.Range("E34") "  -1.31%  "

Set-TextValue $ws.Range("D35") "1.275"
Set-TextValue $ws.Range("E35") "  +1.33%  "

Set-TextValue $ws.Range("D36") "0.7257"
Set-TextValue $ws.Range("E36") "  -1.76%  "

Set-TextValue $ws.Range("D37") "2.710"
Set-TextValue $ws.Range("E37") "  -0.46%  "

Set-TextValue $ws.Range("D38") "0.01920"
Set-TextValue $ws.Range("E38") "  -0.09%  "

Set-TextValue $ws.Range("D39") "2.767"
Set-TextValue $ws.Range("E39") "  -0.91%  "

Set-TextValue $ws.Range("D40") "0.4435"
Set-TextValue $ws.Range("E40") "  -0.40%  "

Set-TextValue $ws.Range("D41") "71.94"
Set-TextValue $ws.Range("E41") "  -0.98%  "

Set-TextValue $ws.Range("D42") "5.916"
Set-TextValue $ws.Range("E42") "  -4.95%  "

Set-TextValue $ws.Range("D43") "0.8465"
Set-TextValue $ws.Range("E43") "  +1.09%  "

Set-TextValue $ws.Range("D44") "0.9990"
Set-TextValue $ws.Range("E44") "  +0.03%  "

Set-TextValue $ws.Range("D45") "1.872"
Set-TextValue $ws.Range("E45") "  -1.66%  "

Set-TextValue $ws.Range("D46") "101.09"
Set-TextValue $ws.Range("E46") "  +0.58%  "

Set-TextValue $ws.Range("D47") "7.560"
Set-TextValue $ws.Range("E47") "  -1.40%  "

Set-TextValue $ws.Range("D48") "9.760"
Set-TextValue $ws.Range("E48") "  -0.89%  "

Set-TextValue $ws.Range("D49") "990.15"
Set-TextValue $ws.Range("E49") "  +5.52%  "

Set-TextValue $ws.Range("D50") "2.036.53"
Set-TextValue $ws.Range("E50") "  -0.79%  "

Set-TextValue $ws.Range("D51") "35.89"
Set-TextValue $ws.Range("E51") "  -1.21%  "
